# Backend Conectando a MongoDB
# Extend Table1 with two new columns (MongoDBCollection, PseudoSchema) and
# populate the relevant rows with the Mongo collection each endpoint uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Grow the table to include the two new columns (J:K) ---------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K23"))

# --- Header row (string order matters for shared-string indices) -------
$ws.Range("J1").Value = "MongoDBCollection"
$ws.Range("J4").Value = "productos"
$ws.Range("K1").Value = "PseudoSchema"
$ws.Range("K4").Value = "sku, descripcion, descripcioncorta, metodoenvio, stock, urlfoto, precio, categorias, historicoprecios, estado"
$ws.Range("K4").WrapText = $true

$ws.Range("J5").Value = "productos"
$ws.Range("J6").Value = "productos"
$ws.Range("J7").Value = "productos"
$ws.Range("J8").Value = "productos"

$ws.Range("J9").Value = "ordenes"
$ws.Range("J10").Value = "ordenes"
$ws.Range("J11").Value = "ordenes"

# --- Column width for the new PseudoSchema column -----------------------
$ws.Columns.Item(11).ColumnWidth = 26.666666666666668

# --- Row 4 grows tall to show the wrapped pseudo-schema text ------------
$ws.Rows.Item(4).RowHeight = 85

# --- Selection / scroll position matches the author's saved view --------
$ws.Range("J13").Select()
